$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the date/time value on row 80 (A80)
$ws.Range("A80").Value = 45453.2916666667

# Copy row 80's formatting down to row 81 before filling in values
$ws.Range("A80").Copy()
$ws.Range("A81").PasteSpecial(-4122)

# Append a new data row (row 81)
$ws.Range("A81").Value = 45455.6031712963
$ws.Range("B81").Value = 600
$ws.Range("C81").Value = 6.21999979019165
$ws.Range("D81").Value = 6.21999979019165
$ws.Range("E81").Value = 6.21999979019165
$ws.Range("F81").Value = 6.21999979019165

# Store the adj_close column as text (matches the source data, which keeps
# this column as a shared string rather than a number) without altering
# the cell's number format/style.
$ws.Range("G81").Formula = '="6.21999979019165"'
$ws.Range("G81").Copy()
$ws.Range("G81").PasteSpecial(-4163)

$ws.Range("H81").Value = "PAL.MI"
